# Update NATMI ligand/receptor TPM-derived expression and specificity
# values (columns G,H,M,N,O,P,Q,R,S,T) for data rows 2-7 to reflect the
# recomputed TPM-based statistics.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 0.001364
$ws.Range("H2").Value = 0.004092
$ws.Range("M2").Value = 17.322719
$ws.Range("N2").Value = 51.96815700000001
$ws.Range("O2").Value = 0.03706849169984819
$ws.Range("P2").Value = 0.03706849169984818
$ws.Range("Q2").Value = 0.023628188716
$ws.Range("R2").Value = 0.212653698444
$ws.Range("S2").Value = 0.03706849169984819
$ws.Range("T2").Value = 0.03706849169984818
$ws.Range("G3").Value = 0.001364
$ws.Range("H3").Value = 0.004092
$ws.Range("O3").Value = 0.01800769446635106
$ws.Range("P3").Value = 0.01800769446635106
$ws.Range("Q3").Value = 0.01147846010666667
$ws.Range("R3").Value = 0.10330614096
$ws.Range("S3").Value = 0.01800769446635106
$ws.Range("T3").Value = 0.01800769446635106
$ws.Range("G4").Value = 0.001364
$ws.Range("H4").Value = 0.004092
$ws.Range("M4").Value = 115.6233063333333
$ws.Range("N4").Value = 346.869919
$ws.Range("O4").Value = 0.2474196788117483
$ws.Range("P4").Value = 0.2474196788117483
$ws.Range("Q4").Value = 0.1577101898386667
$ws.Range("R4").Value = 1.419391708548
$ws.Range("S4").Value = 0.2474196788117483
$ws.Range("T4").Value = 0.2474196788117483
$ws.Range("G5").Value = 0.001364
$ws.Range("H5").Value = 0.004092
$ws.Range("M5").Value = 3.245906333333334
$ws.Range("N5").Value = 9.737719
$ws.Range("O5").Value = 0.006945841006579355
$ws.Range("P5").Value = 0.006945841006579355
$ws.Range("Q5").Value = 0.004427416238666667
$ws.Range("R5").Value = 0.039846746148
$ws.Range("S5").Value = 0.006945841006579355
$ws.Range("T5").Value = 0.006945841006579355
$ws.Range("G6").Value = 0.001364
$ws.Range("H6").Value = 0.004092
$ws.Range("M6").Value = 214.223699
$ws.Range("N6").Value = 642.6710969999999
$ws.Range("O6").Value = 0.4584124125255553
$ws.Range("P6").Value = 0.4584124125255552
$ws.Range("Q6").Value = 0.292201125436
$ws.Range("R6").Value = 2.629810128924
$ws.Range("S6").Value = 0.4584124125255553
$ws.Range("T6").Value = 0.4584124125255552
$ws.Range("G7").Value = 0.001364
$ws.Range("H7").Value = 0.004092
$ws.Range("M7").Value = 108.4856083333333
$ws.Range("N7").Value = 325.456825
$ws.Range("O7").Value = 0.2321458814899178
$ws.Range("P7").Value = 0.2321458814899178
$ws.Range("Q7").Value = 0.1479743697666667
$ws.Range("R7").Value = 1.3317693279
$ws.Range("S7").Value = 0.2321458814899178
$ws.Range("T7").Value = 0.2321458814899178
